$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ntf3"
$ws.Range("C2").Value = "Ntrk2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.630712666666667
$ws.Range("H2").Value = 7.892138000000001
$ws.Range("I2").Value = 0.3947434022685045
$ws.Range("J2").Value = 0.3947434022685045
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.293221
$ws.Range("N2").Value = 0.879663
$ws.Range("O2").Value = 0.02611983441994871
$ws.Range("P2").Value = 0.02611983441994871
$ws.Range("Q2").Value = 0.7713801988326668
$ws.Range("R2").Value = 6.942421789494
$ws.Range("S2").Value = 0.01031063230562054
$ws.Range("T2").Value = 0.01031063230562054

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ntf3"
$ws.Range("C3").Value = "Ntrk2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.630712666666667
$ws.Range("H3").Value = 7.892138000000001
$ws.Range("I3").Value = 0.3947434022685045
$ws.Range("J3").Value = 0.3947434022685045
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 10.74310933333333
$ws.Range("N3").Value = 32.229328
$ws.Range("O3").Value = 0.9569854715114954
$ws.Range("P3").Value = 0.9569854715114955
$ws.Range("Q3").Value = 28.26203380258489
$ws.Range("R3").Value = 254.358304223264
$ws.Range("S3").Value = 0.3777637009459767
$ws.Range("T3").Value = 0.3777637009459767

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ntf3"
$ws.Range("C4").Value = "Ntrk2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.630712666666667
$ws.Range("H4").Value = 7.892138000000001
$ws.Range("I4").Value = 0.3947434022685045
$ws.Range("J4").Value = 0.3947434022685045
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.183427
$ws.Range("N4").Value = 0.550281
$ws.Range("O4").Value = 0.01633949433413
$ws.Range("P4").Value = 0.01633949433413
$ws.Range("Q4").Value = 0.4825437323086668
$ws.Range("R4").Value = 4.342893590778001
$ws.Range("S4").Value = 0.006449907584801428
$ws.Range("T4").Value = 0.006449907584801429

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Ntf3"
$ws.Range("C5").Value = "Ntrk2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.630712666666667
$ws.Range("H5").Value = 7.892138000000001
$ws.Range("I5").Value = 0.3947434022685045
$ws.Range("J5").Value = 0.3947434022685045
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.006232666666666667
$ws.Range("N5").Value = 0.018698
$ws.Range("O5").Value = 0.0005551997344257983
$ws.Range("P5").Value = 0.0005551997344257983
$ws.Range("Q5").Value = 0.01639635514711111
$ws.Range("R5").Value = 0.147567196324
$ws.Range("S5").Value = 0.0005551997344257983
$ws.Range("T5").Value = 0.0002191614321058098

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ntf3"
$ws.Range("C6").Value = "Ntrk2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.99311
$ws.Range("H6").Value = 11.97933
$ws.Range("I6").Value = 0.5991736942634763
$ws.Range("J6").Value = 0.5991736942634763
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.293221
$ws.Range("N6").Value = 0.879663
$ws.Range("O6").Value = 0.02611983441994871
$ws.Range("P6").Value = 0.02611983441994871
$ws.Range("Q6").Value = 1.17086370731
$ws.Range("R6").Value = 10.53777336579
$ws.Range("S6").Value = 0.01565031768295097
$ws.Range("T6").Value = 0.01565031768295097

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ntf3"
$ws.Range("C7").Value = "Ntrk2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.99311
$ws.Range("H7").Value = 11.97933
$ws.Range("I7").Value = 0.5991736942634763
$ws.Range("J7").Value = 0.5991736942634763
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 10.74310933333333
$ws.Range("N7").Value = 32.229328
$ws.Range("O7").Value = 0.9569854715114954
$ws.Range("P7").Value = 0.9569854715114955
$ws.Range("Q7").Value = 42.89841731002666
$ws.Range("R7").Value = 386.08575579024
$ws.Range("S7").Value = 0.5734005203220175
$ws.Range("T7").Value = 0.5734005203220175

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Ntf3"
$ws.Range("C8").Value = "Ntrk2"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.99311
$ws.Range("H8").Value = 11.97933
$ws.Range("I8").Value = 0.5991736942634763
$ws.Range("J8").Value = 0.5991736942634763
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.183427
$ws.Range("N8").Value = 0.550281
$ws.Range("O8").Value = 0.01633949433413
$ws.Range("P8").Value = 0.01633949433413
$ws.Range("Q8").Value = 0.7324441879700001
$ws.Range("R8").Value = 6.59199769173
$ws.Range("S8").Value = 0.00979019518257781
$ws.Range("T8").Value = 0.009790195182577814

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Ntf3"
$ws.Range("C9").Value = "Ntrk2"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3.99311
$ws.Range("H9").Value = 11.97933
$ws.Range("I9").Value = 0.5991736942634763
$ws.Range("J9").Value = 0.5991736942634763
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.006232666666666667
$ws.Range("N9").Value = 0.018698
$ws.Range("O9").Value = 0.0005551997344257983
$ws.Range("P9").Value = 0.0005551997344257983
$ws.Range("Q9").Value = 0.02488772359333333
$ws.Range("R9").Value = 0.22398951234
$ws.Range("S9").Value = 0.0003326610759300065
$ws.Range("T9").Value = 0.0003326610759300065

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Ntf3"
$ws.Range("C10").Value = "Ntrk2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.04053866666666667
$ws.Range("H10").Value = 0.121616
$ws.Range("I10").Value = 0.00608290346801924
$ws.Range("J10").Value = 0.006082903468019241
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.293221
$ws.Range("N10").Value = 0.879663
$ws.Range("O10").Value = 0.02611983441994871
$ws.Range("P10").Value = 0.02611983441994871
$ws.Range("Q10").Value = 0.01188678837866667
$ws.Range("R10").Value = 0.106981095408
$ws.Range("S10").Value = 0.0001588844313771943
$ws.Range("T10").Value = 0.0001588844313771944

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Ntf3"
$ws.Range("C11").Value = "Ntrk2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.04053866666666667
$ws.Range("H11").Value = 0.121616
$ws.Range("I11").Value = 0.00608290346801924
$ws.Range("J11").Value = 0.006082903468019241
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 10.74310933333333
$ws.Range("N11").Value = 32.229328
$ws.Range("O11").Value = 0.9569854715114954
$ws.Range("P11").Value = 0.9569854715114955
$ws.Range("Q11").Value = 0.4355113282275554
$ws.Range("R11").Value = 3.919601954047999
$ws.Range("S11").Value = 0.005821250243501303
$ws.Range("T11").Value = 0.005821250243501304

# Row 12
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Ntf3"
$ws.Range("C12").Value = "Ntrk2"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.04053866666666667
$ws.Range("H12").Value = 0.121616
$ws.Range("I12").Value = 0.00608290346801924
$ws.Range("J12").Value = 0.006082903468019241
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.183427
$ws.Range("N12").Value = 0.550281
$ws.Range("O12").Value = 0.01633949433413
$ws.Range("P12").Value = 0.01633949433413
$ws.Range("Q12").Value = 0.007435886010666667
$ws.Range("R12").Value = 0.06692297409600001
$ws.Range("S12").Value = 0.00009939156675076009
$ws.Range("T12").Value = 0.00009939156675076013

# Row 13
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Ntf3"
$ws.Range("C13").Value = "Ntrk2"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.04053866666666667
$ws.Range("H13").Value = 0.121616
$ws.Range("I13").Value = 0.00608290346801924
$ws.Range("J13").Value = 0.006082903468019241
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.006232666666666667
$ws.Range("N13").Value = 0.018698
$ws.Range("O13").Value = 0.0005551997344257983
$ws.Range("P13").Value = 0.0005551997344257983
$ws.Range("Q13").Value = 0.0002526639964444445
$ws.Range("R13").Value = 0.002273975968
$ws.Range("S13").Value = 0.00000337722638998205
$ws.Range("T13").Value = 0.00000337722638998205
